# Applies the 31-Dec-2022 03:00 UTC symbol-list refresh to the crypto tracker sheet.
# G (Hora) increments 2 -> 3 for every data row; several rows get updated Price (D)
# values, and rows 10-20 shift up one rank (Coin name/Link/Volume columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: row number -> column letter -> new value.
$updates = @{
    2 = @{ "D" = "246.22"; "G" = "3" }
    3 = @{ "D" = "25.38"; "G" = "3" }
    4 = @{ "D" = "5.157"; "G" = "3" }
    5 = @{ "G" = "3" }
    6 = @{ "D" = "6.529"; "G" = "3" }
    7 = @{ "D" = "3.014"; "G" = "3" }
    8 = @{ "D" = "0.8173"; "G" = "3" }
    9 = @{ "D" = "0.8412"; "G" = "3" }
    10 = @{ "B" = "One"; "C" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; "D" = "0.009780"; "E" = "9OneONEBestin24h"; "G" = "3" }
    11 = @{ "B" = "WazirX"; "C" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; "D" = "0.1341"; "E" = "10WazirXWRX"; "G" = "3" }
    12 = @{ "D" = "0.02841"; "G" = "3" }
    13 = @{ "D" = "0.09394"; "G" = "3" }
    14 = @{ "D" = "0.001525"; "G" = "3" }
    15 = @{ "B" = "TigerCash"; "C" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; "D" = "0.006192"; "E" = "14TigerCashTCH"; "G" = "3" }
    16 = @{ "B" = "LEO"; "C" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; "D" = "3.512"; "E" = "15LEOLEO"; "G" = "3" }
    17 = @{ "B" = "BTSEToken"; "C" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; "D" = "2.082"; "E" = "16BTSETokenBTSE"; "G" = "3" }
    18 = @{ "B" = "BitpandaEcosystemToken"; "C" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; "D" = "0.3179"; "E" = "17BitpandaEcosystemTokenBEST"; "G" = "3" }
    19 = @{ "B" = "MandalaExchangeToken"; "C" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; "D" = "0.06954"; "E" = "18MandalaExchangeTokenMDX"; "G" = "3" }
    20 = @{ "B" = "LiechtensteinCryptoassetsExchange"; "C" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; "D" = "0.03166"; "E" = "19LiechtensteinCryptoassetsExchangeLCX"; "G" = "3" }
    21 = @{ "G" = "3" }
    22 = @{ "D" = "3.737"; "G" = "3" }
    23 = @{ "D" = "0.04717"; "G" = "3" }
    24 = @{ "G" = "3" }
    25 = @{ "D" = "0.001243"; "G" = "3" }
    26 = @{ "D" = "0.004273"; "G" = "3" }
    27 = @{ "E" = "26NitroExNTX"; "G" = "3" }
    28 = @{ "E" = "27UpBotsUBXTWorstin24h"; "G" = "3" }
    29 = @{ "G" = "3" }
    30 = @{ "G" = "3" }
    31 = @{ "G" = "3" }
    32 = @{ "G" = "3" }
    33 = @{ "G" = "3" }
    34 = @{ "G" = "3" }
    35 = @{ "G" = "3" }
    36 = @{ "G" = "3" }
    37 = @{ "G" = "3" }
    38 = @{ "G" = "3" }
    39 = @{ "G" = "3" }
    40 = @{ "G" = "3" }
    41 = @{ "D" = "0.006209"; "G" = "3" }
    42 = @{ "D" = "0.1053"; "G" = "3" }
    43 = @{ "G" = "3" }
    44 = @{ "G" = "3" }
    45 = @{ "D" = "0.00005303"; "G" = "3" }
    46 = @{ "G" = "3" }
    47 = @{ "G" = "3" }
    48 = @{ "G" = "3" }
    49 = @{ "G" = "3" }
    50 = @{ "G" = "3" }
    51 = @{ "G" = "3" }
}

# Columns D (Price) and G (Hora) hold digit-only text ("246.22", "3", ...).
# Excel auto-detects a bare numeric-looking assignment as a Number, so force
# the cell to Text format first to keep it an inline/shared string like the original.
$textColumns = @("D", "G")

foreach ($rowNum in $updates.Keys) {
    $rowChanges = $updates[$rowNum]
    foreach ($col in $rowChanges.Keys) {
        $cell = $ws.Range("$col$rowNum")
        if ($textColumns -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowChanges[$col]
    }
}
